# Burndown Chart 2 Update
# Update the "Tider" (Dag 1) actual-progress column (G) with the new
# burndown figures, and leave the selection where the user last worked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$ws.Range("G2").Value = 5
$ws.Range("G3").Value = 3
$ws.Range("G6").Value = 8

$ws.Range("G12").Select()
